# Correct for res/com sf bug.
# This updates the base input numbers on SummaryTable_County_wTotals
# (resunits_new, resunits_change, density_hhemp_new, density_hhemp_change,
# plus two tiny Acres floating point corrections) and lets UnitsDensity's
# formulas that reference this sheet recalc automatically.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("UnitsDensity")
$ws2 = $wb.Worksheets.Item("SummaryTable_County_wTotals")

# ---- Weber (row 2) ----
$ws2.Range("I2").Value = 46100
$ws2.Range("J2").Value = 36000
$ws2.Range("L2").Value = 26
$ws2.Range("M2").Value = 16

# ---- Davis (row 3) ----
$ws2.Range("I3").Value = 69000
$ws2.Range("J3").Value = 58700
$ws2.Range("L3").Value = 35
$ws2.Range("M3").Value = 26

# ---- Salt Lake (row 4) ----
$ws2.Range("D4").Value = 18332.353366010699
$ws2.Range("I4").Value = 234400
$ws2.Range("J4").Value = 189400
$ws2.Range("L4").Value = 45
$ws2.Range("M4").Value = 26

# ---- Utah (row 5) ----
$ws2.Range("I5").Value = 88100
$ws2.Range("J5").Value = 69500
$ws2.Range("L5").Value = 41
$ws2.Range("M5").Value = 25

# ---- Total (row 6) ----
$ws2.Range("D6").Value = 36374.373105956198
$ws2.Range("I6").Value = 437500
$ws2.Range("J6").Value = 353500
$ws2.Range("L6").Value = 40
$ws2.Range("M6").Value = 24

# ---- Selection / active sheet bookkeeping ----
# UnitsDensity keeps its own cursor position (no longer the active tab)
$ws1.Range("D30").Select() | Out-Null

# SummaryTable_County_wTotals becomes the active sheet/tab
$ws2.Activate() | Out-Null
$ws2.Range("A1:M6").Select() | Out-Null
